# Update the "why" (column E) entry for the CIVN2020 congress row (row 2):
# replace the old conference-website link with the DOI link for the
# "Behav Proc in press" publication, and give that cell the new
# dark-grey font style used for it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "http://doi.org/10.17605/OSF.IO/5BWNX"
$ws.Range("E2").Font.Color = 3355443

# Match the printer/page setup metadata written by Excel on save.
$ws.PageSetup.Orientation = 1

# Reflect where the user's selection ended up after the edit.
[void]$ws.Range("E2").Select()
